# Auto-generated edit script: apply updated crypto price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'257.64"
$ws.Range('E2').Value = "'1.24%"
$ws.Range('D3').Value = "'26.96"
$ws.Range('E3').Value = "'-4.02%"
$ws.Range('D4').Value = "'4.750"
$ws.Range('E4').Value = "'-10.72%"
$ws.Range('D5').Value = "'0.05976"
$ws.Range('E5').Value = "'2.17%"
$ws.Range('E6').Value = "'-0.58%"
$ws.Range('D7').Value = "'0.8706"
$ws.Range('E7').Value = "'0.49%"
$ws.Range('D8').Value = "'0.9489"
$ws.Range('E8').Value = "'4.13%"
$ws.Range('D9').Value = "'0.1408"
$ws.Range('E9').Value = "'-1.43%"
$ws.Range('D10').Value = "'0.03609"
$ws.Range('E10').Value = "'4.76%"
$ws.Range('D11').Value = "'0.07188"
$ws.Range('E11').Value = "'0.24%"
$ws.Range('D12').Value = "'0.03172"
$ws.Range('E12').Value = "'-0.30%"
$ws.Range('D13').Value = "'0.09243"
$ws.Range('E13').Value = "'0.24%"
$ws.Range('D14').Value = "'0.001536"
$ws.Range('E14').Value = "'-0.87%"
$ws.Range('D15').Value = "'0.0006072"
$ws.Range('E15').Value = "'-0.05%"
$ws.Range('D16').Value = "'0.005877"
$ws.Range('E16').Value = "'0.12%"
$ws.Range('D17').Value = "'3.485"
$ws.Range('E17').Value = "'-0.42%"
$ws.Range('D18').Value = "'3.195"
$ws.Range('E18').Value = "'-1.01%"
$ws.Range('E19').Value = "'0.77%"
$ws.Range('D20').Value = "'0.3107"
$ws.Range('E20').Value = "'-1.97%"
$ws.Range('D21').Value = "'0.1308"
$ws.Range('E21').Value = "'-0.52%"
$ws.Range('D22').Value = "'3.799"
$ws.Range('E22').Value = "'6.53%"
$ws.Range('D23').Value = "'0.04223"
$ws.Range('E23').Value = "'1.60%"
$ws.Range('E24').Value = "'0.19%"
$ws.Range('D25').Value = "'0.001228"
$ws.Range('E25').Value = "'0.42%"
$ws.Range('D26').Value = "'0.004504"
$ws.Range('E26').Value = "'-10.61%"
$ws.Range('D27').Value = "'0.0001701"
$ws.Range('E27').Value = "'41.69%"
$ws.Range('E28').Value = "'-22.87%"
$ws.Range('D40').Value = "'0.03824"
$ws.Range('E40').Value = "'-0.58%"
$ws.Range('B41').Value = "'BKEXToken"
$ws.Range('C41').Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range('D41').Value = "'0.1103"
$ws.Range('E41').Value = "'0.18%"
$ws.Range('B42').Value = "'KickToken"
$ws.Range('C42').Value = "'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range('D42').Value = "'0.004004"
$ws.Range('E42').Value = "'-29.27%"
$ws.Range('D43').Value = "'0.002301"
$ws.Range('E43').Value = "'-3.34%"
$ws.Range('D44').Value = "'0.01097"
$ws.Range('E44').Value = "'0.59%"
$ws.Range('E45').Value = "'5.05%"
$ws.Range('E46').Value = "'0.12%"
$ws.Range('E47').Value = "'21.56%"
$ws.Range('D48').Value = "'0.002276"
$ws.Range('E48').Value = "'5.50%"
$ws.Range('E49').Value = "'0.12%"
$ws.Range('E50').Value = "'0.12%"
